$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures per the latest data refresh.
# Column D (Price) values that look numeric must be forced to Text so they
# keep being stored as literal strings (matching the sheet's original
# inlineStr/text cell type) instead of being auto-parsed into numbers.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '43.834.86'
$ws.Range('E2').Value = '  +0.28%  '
Set-TextValue 'D3' '2.239.11'
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '270.72'
$ws.Range('E5').Value = '  +4.70%  '
Set-TextValue 'D6' '92.58'
$ws.Range('E6').Value = '  +13.88%  '
Set-TextValue 'D7' '0.628'
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('E8').Value = '  -0.02%  '
Set-TextValue 'D9' '0.625'
$ws.Range('E9').Value = '  +5.47%  '
Set-TextValue 'D10' '46.39'
$ws.Range('E10').Value = '  +8.36%  '
Set-TextValue 'D11' '0.0960'
$ws.Range('E11').Value = '  +4.30%  '
Set-TextValue 'D12' '8.26'
$ws.Range('E12').Value = '  +18.75%  '
$ws.Range('E13').Value = '  +1.27%  '
Set-TextValue 'D14' '2.572.71'
$ws.Range('E14').Value = '  +2.14%  '
Set-TextValue 'D15' '15.00'
$ws.Range('E15').Value = '  +4.90%  '
Set-TextValue 'D16' '2.233.96'
$ws.Range('E16').Value = '  +2.99%  '
Set-TextValue 'D17' '0.801'
$ws.Range('E17').Value = '  +3.11%  '
Set-TextValue 'D18' '43.810.58'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  +2.81%  '
Set-TextValue 'D20' '6.11'
$ws.Range('E20').Value = '  +3.17%  '
Set-TextValue 'D21' '70.85'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('E22').Value = '  -2.25%  '
Set-TextValue 'D23' '234.28'
$ws.Range('E23').Value = '  +1.75%  '
Set-TextValue 'D24' '9.09'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('E25').Value = '  +0.02%  '
Set-TextValue 'D26' '11.42'
$ws.Range('E26').Value = '  +7.82%  '
$ws.Range('E27').Value = '  +13.05%  '
Set-TextValue 'D28' '41.80'
$ws.Range('E28').Value = '  +0.98%  '
Set-TextValue 'D29' '3.55'
$ws.Range('E29').Value = '  +5.31%  '
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  +5.46%  '
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('E34').Value = '  +4.89%  '
$ws.Range('E35').Value = '  +0.92%  '
$ws.Range('E36').Value = '  +1.53%  '
Set-TextValue 'D37' '0.0351'
$ws.Range('E37').Value = '  -0.47%  '
Set-TextValue 'D38' '4.33'
$ws.Range('E38').Value = '  -3.60%  '
Set-TextValue 'D39' '3.57'
$ws.Range('E39').Value = '  +26.22%  '
Set-TextValue 'D40' '0.232'
$ws.Range('E40').Value = '  +17.12%  '
Set-TextValue 'D41' '12.82'
$ws.Range('E41').Value = '  -2.33%  '
$ws.Range('E42').Value = '  +4.80%  '
Set-TextValue 'D43' '63.64'
$ws.Range('E43').Value = '  +1.80%  '
Set-TextValue 'D44' '5.38'
$ws.Range('E44').Value = '  -1.08%  '
Set-TextValue 'D45' '0.0997'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D46' '8.35'
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '100.02'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('E48').Value = '  +4.52%  '
$ws.Range('E49').Value = '  +2.12%  '
Set-TextValue 'D50' '0.448'
$ws.Range('E50').Value = '  +2.50%  '
Set-TextValue 'D51' '2.460.10'
$ws.Range('E51').Value = '  +2.21%  '
